$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header values in row 1 (columns P and Q), matching style of existing header cells (O1)
$ws.Range("P1").Value = 14
$ws.Range("O1").Copy()
$ws.Range("P1").PasteSpecial(-4122)

$ws.Range("Q1").Value = 15
$ws.Range("O1").Copy()
$ws.Range("Q1").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# For each data row (2-25):
#  - swap I <-> K values
#  - swap M <-> O values
#  - add new P, Q columns both = 2
for ($r = 2; $r -le 25; $r++) {
    $iVal = $ws.Cells.Item($r, 9).Value2   # column I
    $kVal = $ws.Cells.Item($r, 11).Value2  # column K
    $mVal = $ws.Cells.Item($r, 13).Value2  # column M
    $oVal = $ws.Cells.Item($r, 15).Value2  # column O

    $ws.Cells.Item($r, 9).Value = $kVal    # I = old K
    $ws.Cells.Item($r, 11).Value = $iVal   # K = old I
    $ws.Cells.Item($r, 13).Value = $oVal   # M = old O
    $ws.Cells.Item($r, 15).Value = $mVal   # O = old M

    $ws.Cells.Item($r, 16).Value = 2       # column P
    $ws.Cells.Item($r, 17).Value = 2       # column Q
}

$ws.Range("A1").Select()
